$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.725.34'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '2.469.29'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.551'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.67%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +3.32%  '
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.09%  '
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '2.849.45'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.78'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").Value = '2.468.03'
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("E17").Value = '  +3.96%  '
$ws.Range("D18").Value = '41.689.86'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("E20").Value = '  +2.39%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.40%  '
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("E32").Value = '  +1.43%  '
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0766'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.30%  '
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("E40").Value = '  -2.38%  '
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '1.971.02'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("D48").Value = '2.702.19'
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.56%  '
